# 2016_19.xlsx — "Mise à jour du TODO. Correction de bug dans la validation
# du panier. Séparation des états de factures des producteurs"
#
# The producer banner changes from "La ferme de Robert" to "La ferme de
# Maurice", the "Pain complet" line item is dropped entirely (the other
# items shift up one row), and the unit prices for the remaining items are
# reset to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the merged producer banner (A10:F11) to the new farm name.
$ws.Range("A10").Value = "La ferme de Maurice"

# 2) Remove the "Pain complet" line (row 12). Excel shifts the rows below
#    it up by one, which also re-points the SUBTOTAL/shared-formula
#    references, the mergeCells list and the sheet dimension automatically.
$ws.Rows.Item(12).Delete()

# 3) After the shift the remaining rows (now 12-15: Pomme de terre, Radis,
#    Salade, Tomates grappe) keep their original quantities, but every
#    unit price (column E) is reset to 1.
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 1
